$wb = $excel.ActiveWorkbook

# --- 1. Append the new log row (row 43) to the "Logs" sheet ---
$logs = $wb.Worksheets.Item("Logs")

$row = 43
$logs.Cells.Item($row, 1).Value = "Kun je 10 dozen schroeven bestellen?"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Testmail #11: Kun je 10 dozen schroeven bestellen?"
$logs.Cells.Item($row, 4).Value = "Bestelling / Levering"
$logs.Cells.Item($row, 5).Value = "Geachte klant,`nBedankt voor uw e-mail. Helaas kan ik geen bestellingen plaatsen, maar ik kan u doorverwijzen naar het bestelteam binnen ons bedrijf. Graag ontvang ik de contactgegevens van uw bedrijf, zodat ik de juiste persoon met u in contact kan brengen.`nIk zie uw reactie graag tegemoet.`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$logs.Cells.Item($row, 6).Value = "2025-06-26 23:33:11"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Nee"
$logs.Cells.Item($row, 9).Value = "Ja"

# Keep the row at the sheet's standard height (writing multi-line text can
# trigger an auto row-height bump; the source workbook has no per-row height
# override, so re-fit it back down to the default).
$logs.Rows.Item($row).AutoFit()

# --- 2. Extend the conditional-formatting ranges so they cover the new row ---
$dFc = $logs.Range("D2:D42").FormatConditions.Item(1)
$dFc.ModifyAppliesToRange($logs.Range("D2:D43"))

$gFc = $logs.Range("G2:G42").FormatConditions.Item(1)
$gFc.ModifyAppliesToRange($logs.Range("G2:G43"))

$hFc = $logs.Range("H2:H42").FormatConditions.Item(1)
$hFc.ModifyAppliesToRange($logs.Range("H2:H43"))

$iFc = $logs.Range("I2:I42").FormatConditions.Item(1)
$iFc.ModifyAppliesToRange($logs.Range("I2:I43"))

# --- 3. Update the Dashboard summary count for "Bestelling / Levering" ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 19
